$d = $word.ActiveDocument

# --- Change 1: insert a new "Meta description" paragraph right after the title (paragraph 1) ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)

$insPoint = $p2.Range.Duplicate
$insPoint.Collapse(1)

$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of 3 Fruits Win: 10 Lines, a classic slot game by Playson. Play for free and discover its simple gameplay and impressive winning potential.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insPoint.InsertXML($metaXml)

# --- Change 2: remove the trailing bold "Play 3 Fruits Win: 10 Lines for Free | Game Review" paragraph ---
# (search from the end backwards, and skip paragraph 1 which is the real Heading1 title)
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Play 3 Fruits Win: 10 Lines for Free | Game Review`r") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    # fall back: the duplicate title paragraph is always the second-to-last paragraph
    $target = $d.Paragraphs($d.Paragraphs.Count - 1)
}
$target.Range.Delete()

# --- Change 3: rewrite the final (italic) paragraph's text ---
$promptPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Read our review of 3 Fruits Win: 10 Lines*") {
        $promptPara = $p
        break
    }
}
if ($promptPara -eq $null) {
    $promptPara = $d.Paragraphs($d.Paragraphs.Count)
}
$rng = $promptPara.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = 'Prompt: Create a feature image for "3 Fruits Win: 10 Lines" that reflects the fun and excitement of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be surrounded by vibrant fruits and stars, with the game title prominently displayed.'
